$d = $word.ActiveDocument

function Get-ParagraphIndexByPattern($doc, $pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# Change 1: insert a new paragraph before "Login to SSH terminal..."
# with the text about network security / port rules (same Arial
# formatting as the paragraph it precedes).
# ------------------------------------------------------------------
$loginParaIndex = Get-ParagraphIndexByPattern $d "^Login to SSH terminal"

if ($loginParaIndex -ge 1) {
    $loginPara = $d.Paragraphs($loginParaIndex)
    $loginPara.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs($loginParaIndex)
    $newPara.Range.Text = "Ensure network security rules allow port 3000 " + [char]0x2013 + " 10000 for inbound traffic"
}

# ------------------------------------------------------------------
# Change 2: split the "navigate to project Directory" run into
# "navigate" + proofErr(gramStart/gramEnd) + " to project Directory"
# (keeps the same Consolas run formatting, just wraps "navigate" with
# grammar-check proofErr markers).
# ------------------------------------------------------------------
$navParaIndex = Get-ParagraphIndexByPattern $d "navigate to project Directory"

if ($navParaIndex -ge 1) {
    $navPara = $d.Paragraphs($navParaIndex)
    $paraStart = $navPara.Range.Start
    $paraEnd = $navPara.Range.End
    $paraText = $d.Range($paraStart, $paraEnd).Text

    $needle = "navigate to project Directory"
    $relIdx = $paraText.IndexOf($needle)
    $absStart = $paraStart + $relIdx
    $absEnd = $absStart + $needle.Length

    $splitRange = $d.Range($absStart, $absEnd)

    $frag = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="848BBD"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w14:ligatures w14:val="none"/></w:rPr><w:t>navigate</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="848BBD"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> to project Directory</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $splitRange.InsertXML($frag)
}

# ------------------------------------------------------------------
# Change 3: fix "./ generate.sh" -> "./generate.sh" (only the text of
# that single run should change; use InsertXML rather than a plain
# text assignment so the adjoining single-space run is left intact
# and not coalesced into this one).
# ------------------------------------------------------------------
$genParaIndex = Get-ParagraphIndexByPattern $d "\./ generate\.sh"

if ($genParaIndex -ge 1) {
    $genPara = $d.Paragraphs($genParaIndex)
    $paraStart = $genPara.Range.Start
    $paraEnd = $genPara.Range.End
    $paraText = $d.Range($paraStart, $paraEnd).Text

    $needle = "./ generate.sh"
    $relIdx = $paraText.IndexOf($needle)
    $absStart = $paraStart + $relIdx
    $absEnd = $absStart + $needle.Length

    $genRange = $d.Range($absStart, $absEnd)

    $frag2 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="FF7EDB"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w14:ligatures w14:val="none"/></w:rPr><w:t>./generate.sh</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $genRange.InsertXML($frag2)
}

Write-Host "Edits applied"
